# Apply the malta_premier-league_2023-2024 update:
# 1) Swap columns F:V between rows 40/41 and rows 46/47 (A:E unchanged)
# 2) Append 7 new match rows (50-56), reusing the existing column A / column E
#    cell formatting (bold/bordered index style + datetime number format) via
#    Copy/PasteSpecial(xlPasteFormats) from the last existing data row (49) so no
#    new style entries are created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap content (columns F:V) between rows 40 and 41 ---
$ws.Range("F40").Value = 'Balzan'
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 'Birkirkara'
$ws.Range("I40").Value = 3
$ws.Range("J40").Value = 2.8
$ws.Range("K40").Value = '28/10/2023 09:43'
$ws.Range("L40").Value = 3.6
$ws.Range("M40").Value = '29/10/2023 14:51'
$ws.Range("N40").Value = 3.03
$ws.Range("O40").Value = '28/10/2023 09:43'
$ws.Range("P40").Value = 2.56
$ws.Range("Q40").Value = '29/10/2023 14:51'
$ws.Range("R40").Value = 2.34
$ws.Range("S40").Value = '28/10/2023 09:43'
$ws.Range("T40").Value = 2.49
$ws.Range("U40").Value = '29/10/2023 14:51'
$ws.Range("V40").Value = 'https://www.betexplorer.com/football/malta/premier-league/balzan-fc-birkirkara/WnBN3OYq/'

$ws.Range("F41").Value = 'Gudja'
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 'Floriana'
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = 4.82
$ws.Range("K41").Value = '28/10/2023 09:43'
$ws.Range("L41").Value = 9.960000000000001
$ws.Range("M41").Value = '29/10/2023 14:43'
$ws.Range("N41").Value = 3.62
$ws.Range("O41").Value = '28/10/2023 09:43'
$ws.Range("P41").Value = 5.01
$ws.Range("Q41").Value = '29/10/2023 14:43'
$ws.Range("R41").Value = 1.58
$ws.Range("S41").Value = '28/10/2023 09:43'
$ws.Range("T41").Value = 1.3
$ws.Range("U41").Value = '29/10/2023 10:32'
$ws.Range("V41").Value = 'https://www.betexplorer.com/football/malta/premier-league/gudja-floriana/UwPe84QS/'

# --- Swap content (columns F:V) between rows 46 and 47 ---
$ws.Range("F46").Value = 'Santa Lucia'
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 'Sliema'
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 3.87
$ws.Range("K46").Value = '05/11/2023 12:12'
$ws.Range("L46").Value = 4.4
$ws.Range("M46").Value = '05/11/2023 13:58'
$ws.Range("N46").Value = 3.34
$ws.Range("O46").Value = '05/11/2023 12:12'
$ws.Range("P46").Value = 3.53
$ws.Range("Q46").Value = '05/11/2023 13:58'
$ws.Range("R46").Value = 1.94
$ws.Range("S46").Value = '05/11/2023 12:12'
$ws.Range("T46").Value = 1.78
$ws.Range("U46").Value = '05/11/2023 13:58'
$ws.Range("V46").Value = 'https://www.betexplorer.com/football/malta/premier-league/santa-lucia-sliema/YPALOPA9/'

$ws.Range("F47").Value = 'Hibernians'
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 'Gudja'
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 1.47
$ws.Range("K47").Value = '04/11/2023 02:43'
$ws.Range("L47").Value = 1.61
$ws.Range("M47").Value = '05/11/2023 13:53'
$ws.Range("N47").Value = 3.95
$ws.Range("O47").Value = '04/11/2023 02:43'
$ws.Range("P47").Value = 3.46
$ws.Range("Q47").Value = '05/11/2023 13:53'
$ws.Range("R47").Value = 5.49
$ws.Range("S47").Value = '04/11/2023 02:43'
$ws.Range("T47").Value = 6.15
$ws.Range("U47").Value = '05/11/2023 13:53'
$ws.Range("V47").Value = 'https://www.betexplorer.com/football/malta/premier-league/hibernians-gudja/CG9PNqQF/'

# --- Append new rows 50-56 ---
# Row 50
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A50").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E50").PasteSpecial(-4122) | Out-Null
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 'malta'
$ws.Range("C50").Value = 'premier-league'
$ws.Range("D50").Value = '2023-2024'
$ws.Range("E50").Value = 45255.54166666666
$ws.Range("F50").Value = 'Sirens'
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 'Hibernians'
$ws.Range("I50").Value = 2
$ws.Range("J50").Value = 5.16
$ws.Range("K50").Value = '24/11/2023 02:12'
$ws.Range("L50").Value = 4.58
$ws.Range("M50").Value = '25/11/2023 12:56'
$ws.Range("N50").Value = 3.76
$ws.Range("O50").Value = '24/11/2023 02:12'
$ws.Range("P50").Value = 3.57
$ws.Range("Q50").Value = '25/11/2023 12:56'
$ws.Range("R50").Value = 1.52
$ws.Range("S50").Value = '24/11/2023 02:12'
$ws.Range("T50").Value = 1.75
$ws.Range("U50").Value = '25/11/2023 12:56'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/malta/premier-league/sirens-hibernians/vRg4G1Ik/'

# Row 51
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A51").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E51").PasteSpecial(-4122) | Out-Null
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 'malta'
$ws.Range("C51").Value = 'premier-league'
$ws.Range("D51").Value = '2023-2024'
$ws.Range("E51").Value = 45255.58333333334
$ws.Range("F51").Value = 'Sliema'
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 'Floriana'
$ws.Range("I51").Value = 2
$ws.Range("J51").Value = 5
$ws.Range("K51").Value = '25/11/2023 11:13'
$ws.Range("L51").Value = 3.87
$ws.Range("M51").Value = '25/11/2023 13:57'
$ws.Range("N51").Value = 3.76
$ws.Range("O51").Value = '25/11/2023 11:13'
$ws.Range("P51").Value = 3.27
$ws.Range("Q51").Value = '25/11/2023 13:57'
$ws.Range("R51").Value = 1.63
$ws.Range("S51").Value = '25/11/2023 11:13'
$ws.Range("T51").Value = 1.97
$ws.Range("U51").Value = '25/11/2023 13:57'
$ws.Range("V51").Value = 'https://www.betexplorer.com/football/malta/premier-league/sliema-floriana/QNk0Hs3q/'

# Row 52
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A52").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E52").PasteSpecial(-4122) | Out-Null
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 'malta'
$ws.Range("C52").Value = 'premier-league'
$ws.Range("D52").Value = '2023-2024'
$ws.Range("E52").Value = 45255.64583333334
$ws.Range("F52").Value = 'Gudja'
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 'Balzan'
$ws.Range("I52").Value = 2
$ws.Range("J52").Value = 3.58
$ws.Range("K52").Value = '24/11/2023 04:42'
$ws.Range("L52").Value = 4.84
$ws.Range("M52").Value = '25/11/2023 15:28'
$ws.Range("N52").Value = 3.18
$ws.Range("O52").Value = '24/11/2023 04:42'
$ws.Range("P52").Value = 2.81
$ws.Range("Q52").Value = '25/11/2023 15:28'
$ws.Range("R52").Value = 1.95
$ws.Range("S52").Value = '24/11/2023 04:42'
$ws.Range("T52").Value = 1.96
$ws.Range("U52").Value = '25/11/2023 15:28'
$ws.Range("V52").Value = 'https://www.betexplorer.com/football/malta/premier-league/gudja-balzan-fc/Mkq9FLXe/'

# Row 53
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E53").PasteSpecial(-4122) | Out-Null
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 'malta'
$ws.Range("C53").Value = 'premier-league'
$ws.Range("D53").Value = '2023-2024'
$ws.Range("E53").Value = 45255.67708333334
$ws.Range("F53").Value = 'Marsaxlokk'
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 'Hamrun'
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 4.19
$ws.Range("K53").Value = '24/11/2023 04:42'
$ws.Range("L53").Value = 5.55
$ws.Range("M53").Value = '25/11/2023 16:13'
$ws.Range("N53").Value = 3.33
$ws.Range("O53").Value = '24/11/2023 04:42'
$ws.Range("P53").Value = 2.83
$ws.Range("Q53").Value = '25/11/2023 16:13'
$ws.Range("R53").Value = 1.74
$ws.Range("S53").Value = '24/11/2023 04:42'
$ws.Range("T53").Value = 1.86
$ws.Range("U53").Value = '25/11/2023 16:13'
$ws.Range("V53").Value = 'https://www.betexplorer.com/football/malta/premier-league/marsaxlokk-hamrun/0trDEum2/'

# Row 54
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A54").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E54").PasteSpecial(-4122) | Out-Null
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 'malta'
$ws.Range("C54").Value = 'premier-league'
$ws.Range("D54").Value = '2023-2024'
$ws.Range("E54").Value = 45256.45833333334
$ws.Range("F54").Value = 'Mosta'
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 'Gzira'
$ws.Range("I54").Value = 3
$ws.Range("J54").Value = 3.08
$ws.Range("K54").Value = '24/11/2023 23:13'
$ws.Range("L54").Value = 3.13
$ws.Range("M54").Value = '26/11/2023 09:47'
$ws.Range("N54").Value = 3.26
$ws.Range("O54").Value = '24/11/2023 23:13'
$ws.Range("P54").Value = 3.3
$ws.Range("Q54").Value = '26/11/2023 10:57'
$ws.Range("R54").Value = 2.07
$ws.Range("S54").Value = '24/11/2023 23:13'
$ws.Range("T54").Value = 2.23
$ws.Range("U54").Value = '26/11/2023 10:57'
$ws.Range("V54").Value = 'https://www.betexplorer.com/football/malta/premier-league/mosta-fc-gzira/SnWfwPP2/'

# Row 55
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A55").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E55").PasteSpecial(-4122) | Out-Null
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 'malta'
$ws.Range("C55").Value = 'premier-league'
$ws.Range("D55").Value = '2023-2024'
$ws.Range("E55").Value = 45256.58333333334
$ws.Range("F55").Value = 'Naxxar'
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 'Birkirkara'
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5.18
$ws.Range("K55").Value = '25/11/2023 02:12'
$ws.Range("L55").Value = 6
$ws.Range("M55").Value = '26/11/2023 09:47'
$ws.Range("N55").Value = 3.91
$ws.Range("O55").Value = '25/11/2023 02:12'
$ws.Range("P55").Value = 4.18
$ws.Range("Q55").Value = '26/11/2023 12:05'
$ws.Range("R55").Value = 1.5
$ws.Range("S55").Value = '25/11/2023 02:12'
$ws.Range("T55").Value = 1.49
$ws.Range("U55").Value = '26/11/2023 09:47'
$ws.Range("V55").Value = 'https://www.betexplorer.com/football/malta/premier-league/naxxar-lions-birkirkara/6wVbxqu9/'

# Row 56
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E56").PasteSpecial(-4122) | Out-Null
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 'malta'
$ws.Range("C56").Value = 'premier-league'
$ws.Range("D56").Value = '2023-2024'
$ws.Range("E56").Value = 45256.67708333334
$ws.Range("F56").Value = 'Valletta'
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 'Santa Lucia'
$ws.Range("I56").Value = 1
$ws.Range("J56").Value = 1.44
$ws.Range("K56").Value = '25/11/2023 04:43'
$ws.Range("L56").Value = 1.49
$ws.Range("M56").Value = '26/11/2023 16:12'
$ws.Range("N56").Value = 4.18
$ws.Range("O56").Value = '25/11/2023 04:43'
$ws.Range("P56").Value = 3.7
$ws.Range("Q56").Value = '26/11/2023 16:12'
$ws.Range("R56").Value = 5.77
$ws.Range("S56").Value = '25/11/2023 04:43'
$ws.Range("T56").Value = 7.92
$ws.Range("U56").Value = '26/11/2023 16:12'
$ws.Range("V56").Value = 'https://www.betexplorer.com/football/malta/premier-league/valletta-santa-lucia/l2K3y3fF/'

$excel.CutCopyMode = $false
